$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of rows 9-15 (columns A:H) while keeping their formatting/style.
$ws.Range("A9:H15").ClearContents()

# Move the active selection to F18 to match the saved view state.
$ws.Range("F18").Select()
